$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="61.127.12"; Text=$true},
    @{Cell="E2"; Value="  -0.22%  "; Text=$false},
    @{Cell="D3"; Value="2.376.85"; Text=$true},
    @{Cell="E3"; Value="  -0.74%  "; Text=$false},
    @{Cell="D4"; Value="0.999"; Text=$true},
    @{Cell="E4"; Value="  -0.11%  "; Text=$false},
    @{Cell="D5"; Value="548.72"; Text=$true},
    @{Cell="E5"; Value="  -0.26%  "; Text=$false},
    @{Cell="D6"; Value="138.23"; Text=$true},
    @{Cell="E6"; Value="  -2.71%  "; Text=$false},
    @{Cell="E7"; Value="  -0.05%  "; Text=$false},
    @{Cell="D8"; Value="0.528"; Text=$true},
    @{Cell="E8"; Value="  -1.90%  "; Text=$false},
    @{Cell="D9"; Value="2.374.51"; Text=$true},
    @{Cell="E9"; Value="  -0.80%  "; Text=$false},
    @{Cell="E10"; Value="  +1.61%  "; Text=$false},
    @{Cell="E11"; Value="  +1.29%  "; Text=$false},
    @{Cell="D12"; Value="5.34"; Text=$true},
    @{Cell="E12"; Value="  +0.89%  "; Text=$false},
    @{Cell="D13"; Value="0.348"; Text=$true},
    @{Cell="E13"; Value="  +0.13%  "; Text=$false},
    @{Cell="D14"; Value="25.06"; Text=$true},
    @{Cell="E14"; Value="  -1.94%  "; Text=$false},
    @{Cell="D15"; Value="2.784.14"; Text=$true},
    @{Cell="E15"; Value="  -1.51%  "; Text=$false},
    @{Cell="D16"; Value="0.0000166"; Text=$true},
    @{Cell="E16"; Value="  -0.22%  "; Text=$false},
    @{Cell="D17"; Value="60.998.62"; Text=$true},
    @{Cell="E17"; Value="  -0.38%  "; Text=$false},
    @{Cell="D18"; Value="2.386.96"; Text=$true},
    @{Cell="E18"; Value="  -0.34%  "; Text=$false},
    @{Cell="D19"; Value="10.80"; Text=$true},
    @{Cell="E19"; Value="  -0.20%  "; Text=$false},
    @{Cell="D20"; Value="4.15"; Text=$true},
    @{Cell="E20"; Value="  +0.13%  "; Text=$false},
    @{Cell="D21"; Value="320.70"; Text=$true},
    @{Cell="E21"; Value="  +0.37%  "; Text=$false},
    @{Cell="D22"; Value="6.71"; Text=$true},
    @{Cell="E22"; Value="  -0.34%  "; Text=$false},
    @{Cell="D24"; Value="64.26"; Text=$true},
    @{Cell="E24"; Value="  +0.87%  "; Text=$false},
    @{Cell="D25"; Value="1.68"; Text=$true},
    @{Cell="E25"; Value="  -12.82%  "; Text=$false},
    @{Cell="D26"; Value="8.37"; Text=$true},
    @{Cell="E26"; Value="  +1.26%  "; Text=$false},
    @{Cell="D27"; Value="0.999"; Text=$true},
    @{Cell="E27"; Value="  -0.20%  "; Text=$false},
    @{Cell="D28"; Value="2.486.37"; Text=$true},
    @{Cell="E28"; Value="  -1.06%  "; Text=$false},
    @{Cell="D29"; Value="8.13"; Text=$true},
    @{Cell="E29"; Value="  +0.13%  "; Text=$false},
    @{Cell="B30"; Value="PEPE"; Text=$false},
    @{Cell="C30"; Value="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; Text=$false},
    @{Cell="D30"; Value="0.0₃0885"; Text=$true},
    @{Cell="E30"; Value="  -5.45%  "; Text=$false},
    @{Cell="B31"; Value="Kaspa"; Text=$false},
    @{Cell="C31"; Value="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; Text=$false},
    @{Cell="D31"; Value="0.150"; Text=$true},
    @{Cell="E31"; Value="  +2.90%  "; Text=$false},
    @{Cell="B32"; Value="Bittensor"; Text=$false},
    @{Cell="C32"; Value="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; Text=$false},
    @{Cell="D32"; Value="504.32"; Text=$true},
    @{Cell="E32"; Value="  -5.23%  "; Text=$false},
    @{Cell="E33"; Value="  -4.60%  "; Text=$false},
    @{Cell="D34"; Value="1.83"; Text=$true},
    @{Cell="E34"; Value="  -1.26%  "; Text=$false},
    @{Cell="E35"; Value="  -4.58%  "; Text=$false},
    @{Cell="D36"; Value="0.998"; Text=$true},
    @{Cell="E36"; Value="  -0.11%  "; Text=$false},
    @{Cell="D37"; Value="4.67"; Text=$true},
    @{Cell="E37"; Value="  -1.32%  "; Text=$false},
    @{Cell="B38"; Value="Stacks"; Text=$false},
    @{Cell="C38"; Value="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; Text=$false},
    @{Cell="D38"; Value="1.88"; Text=$true},
    @{Cell="E38"; Value="  +1.57%  "; Text=$false},
    @{Cell="B39"; Value="PolygonEcosystemToken"; Text=$false},
    @{Cell="C39"; Value="https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; Text=$false},
    @{Cell="D39"; Value="0.378"; Text=$true},
    @{Cell="E39"; Value="  +0.13%  "; Text=$false},
    @{Cell="B40"; Value="RenderToken"; Text=$false},
    @{Cell="C40"; Value="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; Text=$false},
    @{Cell="D40"; Value="5.35"; Text=$true},
    @{Cell="E40"; Value="  -3.81%  "; Text=$false},
    @{Cell="B41"; Value="EthereumClassic"; Text=$false},
    @{Cell="C41"; Value="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; Text=$false},
    @{Cell="D41"; Value="18.57"; Text=$true},
    @{Cell="E41"; Value="  +2.24%  "; Text=$false},
    @{Cell="D42"; Value="146.07"; Text=$true},
    @{Cell="E42"; Value="  +5.34%  "; Text=$false},
    @{Cell="E43"; Value="  -0.06%  "; Text=$false},
    @{Cell="D44"; Value="41.48"; Text=$true},
    @{Cell="E44"; Value="  +2.89%  "; Text=$false},
    @{Cell="D45"; Value="147.56"; Text=$true},
    @{Cell="E45"; Value="  +4.40%  "; Text=$false},
    @{Cell="D46"; Value="3.60"; Text=$true},
    @{Cell="E46"; Value="  -0.87%  "; Text=$false},
    @{Cell="D47"; Value="2.07"; Text=$true},
    @{Cell="E47"; Value="  -4.71%  "; Text=$false},
    @{Cell="D48"; Value="0.0521"; Text=$true},
    @{Cell="E48"; Value="  -0.08%  "; Text=$false},
    @{Cell="B49"; Value="InjectiveProtocol"; Text=$false},
    @{Cell="C49"; Value="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; Text=$false},
    @{Cell="D49"; Value="19.27"; Text=$true},
    @{Cell="E49"; Value="  -5.18%  "; Text=$false},
    @{Cell="B50"; Value="Mantle"; Text=$false},
    @{Cell="C50"; Value="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; Text=$false},
    @{Cell="D50"; Value="0.575"; Text=$true},
    @{Cell="E50"; Value="  -0.61%  "; Text=$false},
    @{Cell="D51"; Value="0.0911"; Text=$true},
    @{Cell="E51"; Value="  +0.35%  "; Text=$false},
)

foreach ($item in $changes) {
    if ($item.Text) {
        $ws.Range($item.Cell).NumberFormat = "@"
    }
    $ws.Range($item.Cell).Value = $item.Value
}
